$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.638.77"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.699.48"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "676.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "3.687.82"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "69.606.19"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "470.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.652"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("D23").Value = "3.845.02"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").Value = "3.689.08"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.945"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("E51").Value = "  +1.87%  "
